# Update column G ("K") values for rows 2-12 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 3
    9  = 2
    10 = 1
    11 = 1
    12 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
